$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing "mahamani140@gmail.com / crossword@123" row
# (old row 3 -> row 4), then insert another new row before the Mobile/Landline
# section (old row 5 -> row 7), leaving a blank gap row behind exactly like the
# original layout had a blank row 4.
$ws.Rows(3).Insert()
$ws.Rows(5).Insert()

# Row 3: new credential pair
$ws.Range("A3").Value = "komsragi@gmail.com"
$ws.Range("B3").Value = "koms@1"

# Row 5: new credential pair
$ws.Range("A5").Value = "komalavalli1998@gmail.com"
$ws.Range("B5").Value = "koms@2905"

# Row 6: keep only A6 styled (matches a stray formatted cell, no value)
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Append the additional phone-number rows after the existing 9876543210 row
$ws.Range("A9").Value = 7490329992
$ws.Range("B9").Value = 4124959590
$ws.Range("A10").Value = 9479200294
$ws.Range("B10").Value = 4132342341
$ws.Range("A11").Value = 6392295295
$ws.Range("B11").Value = 4123399404

$ws.Range("A8:B8").Copy()
$ws.Range("A9:B11").PasteSpecial(-4122)
$ws.Range("A9").Value = 7490329992
$ws.Range("B9").Value = 4124959590
$ws.Range("A10").Value = 9479200294
$ws.Range("B10").Value = 4132342341
$ws.Range("A11").Value = 6392295295
$ws.Range("B11").Value = 4123399404

# Recreate hyperlinks in the same order a user re-adding them by hand would use
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:komsragi@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:mahamani140@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:crossword@123")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:koms@1")
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:komalavalli1998@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:koms@2905")

# Hyperlinks.Add recolors the cell with a brand-new style; restore the shared
# "Hyperlink" look (style used by A2/B2/A3/B3/A4/B4) on every touched cell.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").PasteSpecial(-4122)

# Column A width widened, no longer auto "best fit"
$ws.Columns("A").ColumnWidth = 28.5703125

# Selection moved to D7
$ws.Range("D7").Select()
